$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Backup original formatting of the used range far away, so we can restore
# styles after the merge/unmerge dance (Merge() creates new style entries)
$ws.Range("A1:H23").Copy()
$ws.Range("A1000").PasteSpecial(-4122)

# Drop every merge so we can recreate them in the desired order
$ws.Cells.UnMerge()

# Recreate the merges in the exact target order
$ws.Range("C20:F20").Merge()
$ws.Range("C21:F21").Merge()
$ws.Range("C18:F18").Merge()
$ws.Range("A3:C3").Merge()
$ws.Range("A16:B16").Merge()
$ws.Range("C19:F19").Merge()
$ws.Range("C22:F22").Merge()
$ws.Range("B7:C7").Merge()
$ws.Range("A23:F23").Merge()
$ws.Range("E7:F7").Merge()
$ws.Range("C17:F17").Merge()
$ws.Range("A2:C2").Merge()
$ws.Range("A8:C8").Merge()
$ws.Range("B4:C4").Merge()
$ws.Range("A9:B9").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("D8:F8").Merge()
$ws.Range("E4:F4").Merge()
$ws.Range("D2:E2").Merge()
$ws.Range("C16:F16").Merge()
$ws.Range("A1:F1").Merge()
$ws.Range("A22:B22").Merge()
$ws.Range("E6:F6").Merge()
$ws.Range("B6:C6").Merge()
$ws.Range("E5:F5").Merge()
$ws.Range("A15:F15").Merge()
$ws.Range("D3:F3").Merge()

# Restore the original cell formatting (Merge() rewrites border styles)
$ws.Range("A1000:H1022").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Clean up the scratch backup area
$ws.Range("A1000:H1022").Clear()

# Rename the sheet
$ws.Name = "Sheet1"
